# Target: slide 71 ("Implementing methods parseVariable() and
# parseVariableExpr()"), the "Content Placeholder 2" shape, in the
# bullet:
#     variable = ( varId | paramId) { indexExpr | fieldExpr } .
# The run that holds ") { " gets a space typed in front of the ")"
# and is thereby split into two runs:
#     " ) "   and   "{ "
#
# Note: slide 49 contains a similar-looking line, but with the fix
# already applied (" ) { " as its own run) -- so we match on the
# *unsplit* pattern "paramId) {" (no space between "Id" and ")")
# which uniquely identifies the run that still needs editing.
$pres = $ppt.ActivePresentation

$targetShape = $null

foreach ($sl in $pres.Slides) {
    foreach ($sh in $sl.Shapes) {
        if (-not $sh.HasTextFrame) { continue }
        $tf = $sh.TextFrame
        if (-not $tf.HasText) { continue }
        $txt = $tf.TextRange.Text
        if ($txt.IndexOf("paramId) { indexExpr") -ge 0) {
            $targetShape = $sh
            break
        }
    }
    if ($targetShape -ne $null) { break }
}

if ($targetShape -eq $null) {
    throw "Could not find the shape containing 'paramId) { indexExpr'"
}

$tr = $targetShape.TextFrame.TextRange

# Locate "( varId | paramId) { indexExpr" and split the ") { " run.
$text = $tr.Text
$idx = $text.IndexOf(") { indexExpr")

# Rewrite the lone trailing space of the 4-character run (") { ")
# first, turning it into "{ "; this leaves the run as ") {" + "{ ".
$tail = $tr.Characters($idx + 4, 1)
$tail.Text = "{ "

# Now rewrite the leading 3 characters (") {") of what remains into
# " ) " -- the engine splits this remaining piece off cleanly,
# producing the final pair of runs " ) " followed by "{ ".
$text2 = $tr.Text
$idx2 = $text2.IndexOf(") {")
$head = $tr.Characters($idx2 + 1, 3)
$head.Text = " ) "
